# Add a "Save" column (H) to the s_vals sheet, matching the style used by
# the other header cells (B1:G1) and extending the data row with a 0 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1
# header cell so it reuses the same cell style instead of minting a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and the corresponding data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
